# AN160 TC2: add a new "CE representative details" step to the
# LegalLabelling Gherkin table, right after the "manufacturer details"
# step (row 11), pushing the remaining steps down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LEGALLABELLING")

# Insert a new row at 12 (shifts rows 12-16 down to 13-17, carrying
# their formatting/styles with them, and updates the sheet dimension).
$ws.Rows(12).Insert()

# Populate the newly inserted row with the new Gherkin step.
$ws.Range("C12").Value = "And "
$ws.Range("D12").Value = "the information includes the CE representative details"
